$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the date column as plain text (matches original inlineStr cells)
$ws.Range("A2:A4").NumberFormat = "@"

# Update existing rows 2-4 with new dates (kept as text) and new values
$ws.Range("A2").Value = "05/01/25"
$ws.Range("B2").Value = 1.4

$ws.Range("A3").Value = "06/01/25"
$ws.Range("B3").Value = 1.9

$ws.Range("A4").Value = "07/01/25"
$ws.Range("B4").Value = 2.6

# Remove row 5 entirely, shrinking the used range from A1:B5 to A1:B4
$ws.Rows.Item(5).Delete()
